$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new "students" (team 3) into the previously-empty template rows 13-15
$ws.Range("A13").Value = "Jack"
$ws.Range("B13").Value = 3

$ws.Range("A14").Value = "Jill"
$ws.Range("B14").Value = 3

$ws.Range("A15").Value = "Jane"
$ws.Range("B15").Value = 3

# Move the active selection to B16, matching the author's final cursor position
$ws.Range("B16").Select()
